# Append: 2025-11-19 01:49 JST
# Updates the "取得日時" (acquisition datetime) column A values for rows 2-15
# on the "ランサーズ" sheet from "2025-11-19 01:19:49" to "2025-11-19 01:49:34".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-11-19 01:49:34"
}

$wb.Save()
